# Rename the "SwateTemplateMetadata" sheet to "isa_template" and make it
# the active/selected tab (it was "Sample information" before). Also move
# the remembered selection on "Sample information" from T1 to T29.

$wb = $excel.ActiveWorkbook

# Sheet 2 ("SwateTemplateMetadata") -> "isa_template"
$wsMeta = $wb.Worksheets.Item(2)
$wsMeta.Name = "isa_template"

# Sheet 1 ("Sample information") keeps its own selection, but moves to
# T29 and loses the "selected tab" status.
$wsSample = $wb.Worksheets.Item(1)
$wsSample.Range("T29").Select()

# Activate "isa_template" last so it becomes the active/selected sheet
# (workbook.xml bookViews/workbookView activeTab, and this sheet's
# sheetView gets tabSelected="1").
$wsMeta.Activate()
$wsMeta.Range("F10").Select()
